# Updates cryptos worksheet with latest scraped values (GitHub Actions refresh).
# Numeric-looking text values (column D "Price") are written with a leading
# apostrophe so Excel keeps them as text (matching the source data's formatting,
# e.g. keeping trailing zeros like "1.00" instead of turning them into the number 1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '44.161.76' }
    @{ Cell = 'E2'; Value = '  +2.20%  ' }
    @{ Cell = 'D3'; Value = '2.428.92' }
    @{ Cell = 'E3'; Value = '  +1.94%  ' }
    @{ Cell = 'D4'; Value = '''1.00' }
    @{ Cell = 'D5'; Value = '''307.89' }
    @{ Cell = 'E5'; Value = '  +1.57%  ' }
    @{ Cell = 'D6'; Value = '''100.44' }
    @{ Cell = 'E6'; Value = '  +3.56%  ' }
    @{ Cell = 'E7'; Value = '  +0.59%  ' }
    @{ Cell = 'E8'; Value = '  -0.04%  ' }
    @{ Cell = 'E9'; Value = '  -0.37%  ' }
    @{ Cell = 'D10'; Value = '''35.43' }
    @{ Cell = 'E10'; Value = '  +3.58%  ' }
    @{ Cell = 'E11'; Value = '  +1.33%  ' }
    @{ Cell = 'E12'; Value = '  +2.83%  ' }
    @{ Cell = 'D13'; Value = '''18.71' }
    @{ Cell = 'E13'; Value = '  +1.65%  ' }
    @{ Cell = 'E14'; Value = '  +2.22%  ' }
    @{ Cell = 'D15'; Value = '2.806.85' }
    @{ Cell = 'E15'; Value = '  +1.86%  ' }
    @{ Cell = 'D16'; Value = '2.444.68' }
    @{ Cell = 'E16'; Value = '  +2.93%  ' }
    @{ Cell = 'D17'; Value = '''0.831' }
    @{ Cell = 'E17'; Value = '  +2.73%  ' }
    @{ Cell = 'D18'; Value = '44.122.22' }
    @{ Cell = 'E18'; Value = '  +2.08%  ' }
    @{ Cell = 'D19'; Value = '''12.35' }
    @{ Cell = 'E19'; Value = '  +1.29%  ' }
    @{ Cell = 'D20'; Value = '''6.45' }
    @{ Cell = 'E20'; Value = '  +1.85%  ' }
    @{ Cell = 'D21'; Value = '0.0₃0906' }
    @{ Cell = 'E21'; Value = '  +1.87%  ' }
    @{ Cell = 'D22'; Value = '''68.65' }
    @{ Cell = 'E22'; Value = '  -0.05%  ' }
    @{ Cell = 'D23'; Value = '''2.30' }
    @{ Cell = 'E23'; Value = '  +3.16%  ' }
    @{ Cell = 'D24'; Value = '''240.29' }
    @{ Cell = 'E25'; Value = '  +2.00%  ' }
    @{ Cell = 'E26'; Value = '  -0.04%  ' }
    @{ Cell = 'D27'; Value = '''25.29' }
    @{ Cell = 'E27'; Value = '  +1.90%  ' }
    @{ Cell = 'D29'; Value = '''9.61' }
    @{ Cell = 'E29'; Value = '  +5.25%  ' }
    @{ Cell = 'D30'; Value = '''32.79' }
    @{ Cell = 'E30'; Value = '  +3.99%  ' }
    @{ Cell = 'E31'; Value = '  +15.83%  ' }
    @{ Cell = 'E32'; Value = '  +9.35%  ' }
    @{ Cell = 'E33'; Value = '  +1.71%  ' }
    @{ Cell = 'E34'; Value = '  +0.03%  ' }
    @{ Cell = 'E35'; Value = '  +3.53%  ' }
    @{ Cell = 'E36'; Value = '  +3.38%  ' }
    @{ Cell = 'D37'; Value = '''4.53' }
    @{ Cell = 'E37'; Value = '  +4.80%  ' }
    @{ Cell = 'D38'; Value = '''129.98' }
    @{ Cell = 'E38'; Value = '  +22.79%  ' }
    @{ Cell = 'D39'; Value = '''2.92' }
    @{ Cell = 'E39'; Value = '  +4.66%  ' }
    @{ Cell = 'E40'; Value = '  -0.74%  ' }
    @{ Cell = 'E41'; Value = '  +0.17%  ' }
    @{ Cell = 'D42'; Value = '''21.21' }
    @{ Cell = 'E42'; Value = '  -4.60%  ' }
    @{ Cell = 'E43'; Value = '  +2.80%  ' }
    @{ Cell = 'D44'; Value = '1.962.59' }
    @{ Cell = 'E44'; Value = '  +0.39%  ' }
    @{ Cell = 'D45'; Value = '''2.17' }
    @{ Cell = 'E45'; Value = '  +1.87%  ' }
    @{ Cell = 'D46'; Value = '''2.88' }
    @{ Cell = 'E46'; Value = '  +4.65%  ' }
    @{ Cell = 'D47'; Value = '''9.40' }
    @{ Cell = 'E47'; Value = '  +1.40%  ' }
    @{ Cell = 'E48'; Value = '  +8.54%  ' }
    @{ Cell = 'B49'; Value = 'RocketPoolETH' }
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth' }
    @{ Cell = 'D49'; Value = '2.665.32' }
    @{ Cell = 'E49'; Value = '  +1.77%  ' }
    @{ Cell = 'B50'; Value = 'MultiversX' }
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld' }
    @{ Cell = 'D50'; Value = '''53.35' }
    @{ Cell = 'E50'; Value = '  +1.02%  ' }
    @{ Cell = 'B51'; Value = 'BitcoinSV' }
    @{ Cell = 'C51'; Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv' }
    @{ Cell = 'D51'; Value = '''73.59' }
    @{ Cell = 'E51'; Value = '  +2.33%  ' }
)

foreach ($update in $updates) {
    $ws.Range($update.Cell).Value = $update.Value
}
